$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: add the new "Start" / "End" / "Date" columns (D, E, F) ---
# Copy the existing header formatting (bold font + border + centered/top
# alignment, style index 1) from C1 onto the three new header cells before
# putting their text in, so D1:F1 end up styled the same as A1:C1.
$ws.Range("C1").Copy()
$ws.Range("D1:F1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("D1").Value = "Start"
$ws.Range("E1").Value = "End"
$ws.Range("F1").Value = "Date"

# --- Row 2: existing "Puzzles" / "Router" record gains an ID + Start/End flags ---
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "Puzzles"
$ws.Range("C2").Value = "Router"
$ws.Range("D2").Value = $true
$ws.Range("E2").Value = $false

# --- Row 3: new "Chalkboards" / "Laser" record ---
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "Chalkboards"
$ws.Range("C3").Value = "Laser"
$ws.Range("D3").Value = $true
$ws.Range("E3").Value = $false

# --- Row 4: new "Puzzles" / "Laser" record; ID is the text "12" (not the
# number 12), so type it as a formula and paste the computed text back over
# itself as a value -- that keeps it a genuine string cell instead of
# Excel's normal "numeric-looking text becomes a number" coercion.
$ws.Range("A4").Formula = '="12"'
$ws.Range("A4").Copy()
$ws.Range("A4").PasteSpecial(-4163)  # xlPasteValues

$ws.Range("B4").Value = "Puzzles"
$ws.Range("C4").Value = "Laser"
$ws.Range("D4").Value = $true
$ws.Range("E4").Value = $false
$ws.Range("F4").Value = "04/05/2023, 12:00:48"
